$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateText = "29-11-2025"
$priceText = "The price of gold in India today is " + [char]8377 + "12,982 per gram for 24 karat gold, " + [char]8377 + "11,900 per gram for 22 karat gold and " + [char]8377 + "9,737 per gram for 18 karat gold (also called 999 gold)."

$newRow = 74

$cellA = $ws.Cells.Item($newRow, 1)
$cellA.Value = $dateText
$cellA.Borders.LineStyle = 1

$cellB = $ws.Cells.Item($newRow, 2)
$cellB.Value = $priceText
$cellB.Borders.LineStyle = 1
$cellB.WrapText = $true
